$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Font.Bold = $false
$ws.Rows.Item(2).Font.Name = "Calibri"
$ws.Rows.Item(2).Font.Size = 11
